$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.066970738207348
$ws.Range("D2").Value = 1.070036877672472
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.07865182360091
$ws.Range("I2").Value = 1.05298346427525
$ws.Range("J2").Value = 1.071918021215107
$ws.Range("K2").Value = 1.072737785842959
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.081329935659295
$ws.Range("N2").Value = 1.073440267755948
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.068253724991089
$ws.Range("D3").Value = 1.071091080070438
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.079898811873038
$ws.Range("I3").Value = 1.053392705478733
$ws.Range("J3").Value = 1.072855629750853
$ws.Range("K3").Value = 1.073607684049577
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.082393801327244
$ws.Range("N3").Value = 1.07437920780336
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.069083121575856
$ws.Range("D4").Value = 1.071772494408561
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.080705345145767
$ws.Range("I4").Value = 1.053655888916688
$ws.Range("J4").Value = 1.073461014817454
$ws.Range("K4").Value = 1.074169237580649
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.083081274513331
$ws.Range("N4").Value = 1.074985452586195
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.069431616112836
$ws.Range("D5").Value = 1.072058789488358
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.081044330164534
$ws.Range("I5").Value = 1.053766143787549
$ws.Range("J5").Value = 1.073715207274682
$ws.Range("K5").Value = 1.07440499858878
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.083370070573069
$ws.Range("N5").Value = 1.075240006025868
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.069490119198659
$ws.Range("D6").Value = 1.072106849695111
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.081101242532551
$ws.Range("I6").Value = 1.053784633377352
$ws.Range("J6").Value = 1.073757869110657
$ws.Range("K6").Value = 1.07444456542184
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.083418548006338
$ws.Range("N6").Value = 1.075282728446543
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.069087778894908
$ws.Range("D7").Value = 1.071776320568747
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.080709874997467
$ws.Range("I7").Value = 1.053657363668702
$ws.Range("J7").Value = 1.073464412570366
$ws.Range("K7").Value = 1.074172389070317
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.083085134272775
$ws.Range("N7").Value = 1.074988855164305
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.06740449345947
$ws.Range("D8").Value = 1.070393301703813
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.07907332347059
$ws.Range("I8").Value = 1.053122105969192
$ws.Range("J8").Value = 1.072235163099385
$ws.Range("K8").Value = 1.073032048362555
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.08168966575464
$ws.Range("N8").Value = 1.073757860018079
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.064432192794945
$ws.Range("D9").Value = 1.067950598371751
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.076186701952599
$ws.Range("I9").Value = 1.052166433831301
$ws.Range("J9").Value = 1.070058930303432
$ws.Range("K9").Value = 1.071012354709842
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.07922352729879
$ws.Range("N9").Value = 1.071578536721938
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.062446320985135
$ws.Range("D10").Value = 1.066318198587601
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.074260233941584
$ws.Range("I10").Value = 1.05152085783055
$ws.Range("J10").Value = 1.068601146690097
$ws.Range("K10").Value = 1.069658859786115
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.077574480908121
$ws.Range("N10").Value = 1.070118682888662
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.061585339766081
$ws.Range("D11").Value = 1.065610388511377
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.073425527344039
$ws.Range("I11").Value = 1.051239293263149
$ws.Range("J11").Value = 1.067968226015355
$ws.Range("K11").Value = 1.069071084114488
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.076859216352319
$ws.Range("N11").Value = 1.069484863393964
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.061265365683326
$ws.Range("D12").Value = 1.065347327836501
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.073115396698421
$ws.Range("I12").Value = 1.051134401819372
$ws.Range("J12").Value = 1.067732874136901
$ws.Range("K12").Value = 1.068852499507622
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.076593349423566
$ws.Range("N12").Value = 1.069249177288848
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.06133400884994
$ws.Range("D13").Value = 1.065403762008165
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.073181924604732
$ws.Range("I13").Value = 1.051156915252922
$ws.Range("J13").Value = 1.067783369589104
$ws.Range("K13").Value = 1.068899398390177
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.076650387240388
$ws.Range("N13").Value = 1.069299744450385
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.061558894046398
$ws.Range("D14").Value = 1.065588646879111
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.073399893569844
$ws.Range("I14").Value = 1.051230629150918
$ws.Range("J14").Value = 1.067948777014626
$ws.Range("K14").Value = 1.069053021132331
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.076837243517516
$ws.Range("N14").Value = 1.069465386773423
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.061697430939613
$ws.Range("D15").Value = 1.065702540790619
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.073534180285078
$ws.Range("I15").Value = 1.05127600614758
$ws.Range("J15").Value = 1.068050655850765
$ws.Range("K15").Value = 1.069147638834459
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.076952347121778
$ws.Range("N15").Value = 1.069567410289196
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.062503438240475
$ws.Range("D16").Value = 1.066365152937118
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.074315619150292
$ws.Range("I16").Value = 1.051539501535894
$ws.Range("J16").Value = 1.068643115677784
$ws.Range("K16").Value = 1.06969783240805
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.077621924769208
$ws.Range("N16").Value = 1.070160711477126
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.063008731700154
$ws.Range("D17").Value = 1.066780530354065
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.074805649991491
$ws.Range("I17").Value = 1.051704241813211
$ws.Range("J17").Value = 1.069014295339413
$ws.Range("K17").Value = 1.07004249638976
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.078041605275866
$ws.Range("N17").Value = 1.070532418256449
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.063303356366738
$ws.Range("D18").Value = 1.067022719550443
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.075091425468398
$ws.Range("I18").Value = 1.051800136659307
$ws.Range("J18").Value = 1.06923063522377
$ws.Range("K18").Value = 1.070243368881202
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.078286280639664
$ws.Range("N18").Value = 1.070749065368256
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.063403798173076
$ws.Range("D19").Value = 1.067105284043581
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.075188858933012
$ws.Range("I19").Value = 1.051832801229262
$ws.Range("J19").Value = 1.069304373985883
$ws.Range("K19").Value = 1.070311833441339
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.078369688832313
$ws.Range("N19").Value = 1.070822908847869
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.062954529330654
$ws.Range("D20").Value = 1.066735973972424
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.074753079643388
$ws.Range("I20").Value = 1.051686586959753
$ws.Range("J20").Value = 1.068974488149763
$ws.Range("K20").Value = 1.070005534208099
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.077996589672506
$ws.Range("N20").Value = 1.070492554536023
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.061492675614078
$ws.Range("D21").Value = 1.06553420701406
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.073335709472872
$ws.Range("I21").Value = 1.051208930699293
$ws.Range("J21").Value = 1.067900075767521
$ws.Range("K21").Value = 1.069007790222194
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.076782224155422
$ws.Range("N21").Value = 1.069416616364961
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.060572579979888
$ws.Range("D22").Value = 1.064777747903295
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.072444068110515
$ws.Range("I22").Value = 1.050906839103629
$ws.Range("J22").Value = 1.067223061905572
$ws.Range("K22").Value = 1.068378971765368
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.07601762770704
$ws.Range("N22").Value = 1.068738641065675
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.061060433459628
$ws.Range("D23").Value = 1.06517884367548
$ws.Range("E23").Value = 0.9879432794636464
$ws.Range("F23").Value = 1.072916791108263
$ws.Range("I23").Value = 1.051067151852058
$ws.Range("J23").Value = 1.067582101808917
$ws.Range("K23").Value = 1.06871246312632
$ws.Range("L23").Value = 0.9917760702887611
$ws.Range("M23").Value = 1.076423057825709
$ws.Range("N23").Value = 1.069098190846864
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.062979021370911
$ws.Range("D24").Value = 1.066756107371996
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.074776834079021
$ws.Range("I24").Value = 1.051694565031709
$ws.Range("J24").Value = 1.068992475808626
$ws.Range("K24").Value = 1.070022236334184
$ws.Range("L24").Value = 0.9929938892766441
$ws.Range("M24").Value = 1.078016930647963
$ws.Range("N24").Value = 1.070510567739425
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.065201352768578
$ws.Range("D25").Value = 1.068582778715761
$ws.Range("E25").Value = 0.9912096547607051
$ws.Range("F25").Value = 1.076933312610243
$ws.Range("I25").Value = 1.05241498414026
$ws.Range("J25").Value = 1.070622754995698
$ws.Range("K25").Value = 1.071535723946969
$ws.Range("L25").Value = 0.9944092447426416
$ws.Range("M25").Value = 1.079861945318173
$ws.Range("N25").Value = 1.072143162109939
